$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new Job Posting row (JD_013) at the end of the table (row 14)
$ws.Range("A14").Value = "JD_013"
$ws.Range("B14").Value = "Senior Python Engineer"
$ws.Range("C14").Value = "We are seeking a Software Engineer to build and maintain high-quality software solutions.`nWork with global teams to drive innovation and deliver scalable applications.`nJoin Akkodis and be part of a tech-driven, collaborative environment."
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 7

# Avoid Excel auto-expanding the row height because of the embedded
# newlines in C14 (keep the row looking like the other plain data rows)
$ws.Rows.Item(14).EntireRow.AutoFit()
